# Replace all "OIE" references with "WOAH" across the workbook,
# as described in the commit "OIE replaced with WOAH all Excels".
#
# Note: URLs that happen to contain "oie" (e.g. https://wahis.oie.int/,
# https://www.oie.int/...) are intentionally left untouched, matching
# the source diff.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet 1")
$ws2 = $wb.Worksheets.Item("References")

# -- Sheet 1 --------------------------------------------------------------

$ws1.Range("E5").Value = "Based on official disease reports to the WOAH"

$ws1.Range("E6").Value = "CCHFV is a disease listed in the World Organisation for Animal Health ({ref005:WOAH}) Terrestrial Animal Health Code. The map to the right displays occurrence reported to the {ref001:WOAH-WAHIS} system since 2005."

$ws1.Range("E7").Value = "As described in the WOAH {ref005:Terrestrial Animal Health Code}, the WOAH early warning system includes immediate notifications and follow-up reports on:"

$ws1.Range("E14").Value = "Information on stable situations (disease present or absent in a zone or country) is provided by countries through the WOAH monitoring system, which is a different reporting channel. This information is available in a different spatial and temporal scale, which can be browsed on the map independently from the outbreak notification points."

$ws1.Range("E17").Value = "For more up to date reports, visit the original data source: {ref001:WOAH-WAHIS}."

$nbsp = [char]0x00A0
$ws1.Range("E72").Value = "WOAH-prescribed tests for international trade include ({ref010:WOAH," + $nbsp + "Terrestrial Manual}):"

$ws1.Range("E139").Value = "Geographical distribution data has been kindly provided by the World Organisation of Animal Health (WOAH). {ref001:WOAH-WAHIS} (WOAH World Animal Health Information System) is the original source of these data."

# -- References sheet ------------------------------------------------------

$ws2.Range("C2").Value = "WOAH-WAHIS (WOAH World Animal Health Information System)"

$ws2.Range("C5").Value = "WOAH (World Organisation for Animal Health). Terrestrial Animal Health Code 2021. WOAH, Paris, France"

$ws2.Range("C10").Value = "WOAH (World Organisation for Animal Health), 2018. Crimean-Congo Haemorrhagic Fever. Chapter 3.1.5. WOAH Terrestrial Animal Health Manual, Paris, France"
